$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Replace XNAMATH with a Rorn maths library" row (row 2).
# This shifts all subsequent rows up by one.
$ws.Rows("2").Delete()

# Add the new task at the end of the list (now row 10).
$ws.Range("A10").Value = "Revise, understand and document the view and projection matrix builds"
$ws.Range("B10").Value = 7

# Update the active selection to match the new state (single cell B10).
$ws.Range("B10").Select()
